$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting for new rows 58-64 by copying formats from row 57
$ws.Range("A57").Copy() | Out-Null
$ws.Range("A58:A64").PasteSpecial(-4122) | Out-Null
$ws.Range("B57:E57").Copy() | Out-Null
$ws.Range("B58:E64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(5, 1).Value = 43922
$ws.Cells.Item(5, 2).Value = 0.47890686988830566
$ws.Cells.Item(5, 3).Value = 0.37072709202766418
$ws.Cells.Item(5, 4).Value = 0.60682821273803711
$ws.Cells.Item(5, 5).Value = 0.56078648567199707
$ws.Cells.Item(6, 1).Value = 43952
$ws.Cells.Item(6, 2).Value = 0.5860331654548645
$ws.Cells.Item(6, 3).Value = 0.52606886625289917
$ws.Cells.Item(6, 4).Value = 0.58285558223724365
$ws.Cells.Item(6, 5).Value = 0.52861630916595459
$ws.Cells.Item(7, 1).Value = 43983
$ws.Cells.Item(7, 2).Value = 0.80940514802932739
$ws.Cells.Item(7, 3).Value = 0.73801320791244507
$ws.Cells.Item(7, 4).Value = 0.62259155511856079
$ws.Cells.Item(7, 5).Value = 0.56516802310943604
$ws.Cells.Item(8, 1).Value = 44013
$ws.Cells.Item(8, 2).Value = 0.43972894549369812
$ws.Cells.Item(8, 3).Value = 0.41535386443138123
$ws.Cells.Item(8, 4).Value = 0.59996819496154785
$ws.Cells.Item(8, 5).Value = 0.54142719507217407
$ws.Cells.Item(9, 1).Value = 44044
$ws.Cells.Item(9, 2).Value = 0.72006678581237793
$ws.Cells.Item(9, 3).Value = 0.75376933813095093
$ws.Cells.Item(9, 4).Value = 0.6180424690246582
$ws.Cells.Item(9, 5).Value = 0.56387245655059814
$ws.Cells.Item(10, 1).Value = 44075
$ws.Cells.Item(10, 2).Value = 0.46299266815185547
$ws.Cells.Item(10, 3).Value = 0.36776545643806458
$ws.Cells.Item(10, 4).Value = 0.62297284603118896
$ws.Cells.Item(10, 5).Value = 0.57307511568069458
$ws.Cells.Item(11, 1).Value = 44105
$ws.Cells.Item(11, 2).Value = 0.86100733280181885
$ws.Cells.Item(11, 3).Value = 0.78447830677032471
$ws.Cells.Item(11, 4).Value = 0.63974106311798096
$ws.Cells.Item(11, 5).Value = 0.59626471996307373
$ws.Cells.Item(12, 1).Value = 44136
$ws.Cells.Item(12, 2).Value = 0.44160452485084534
$ws.Cells.Item(12, 3).Value = 0.37524142861366272
$ws.Cells.Item(12, 4).Value = 0.57559704780578613
$ws.Cells.Item(12, 5).Value = 0.53597396612167358
$ws.Cells.Item(13, 1).Value = 44166
$ws.Cells.Item(13, 2).Value = 0.76263660192489624
$ws.Cells.Item(13, 3).Value = 0.74343442916870117
$ws.Cells.Item(13, 4).Value = 0.54548430442810059
$ws.Cells.Item(13, 5).Value = 0.515350341796875
$ws.Cells.Item(14, 1).Value = 44197
$ws.Cells.Item(14, 2).Value = 0.52328032255172729
$ws.Cells.Item(14, 3).Value = 0.45355129241943359
$ws.Cells.Item(14, 4).Value = 0.46770617365837097
$ws.Cells.Item(14, 5).Value = 0.42931428551673889
$ws.Cells.Item(15, 1).Value = 44228
$ws.Cells.Item(15, 2).Value = 0.73694723844528198
$ws.Cells.Item(15, 3).Value = 0.73477494716644287
$ws.Cells.Item(15, 4).Value = 0.52602195739746094
$ws.Cells.Item(15, 5).Value = 0.49583402276039124
$ws.Cells.Item(16, 1).Value = 44256
$ws.Cells.Item(16, 2).Value = 0.23210926353931427
$ws.Cells.Item(16, 3).Value = 0.19539661705493927
$ws.Cells.Item(16, 4).Value = 0.52616488933563232
$ws.Cells.Item(16, 5).Value = 0.49776607751846313
$ws.Cells.Item(17, 1).Value = 44287
$ws.Cells.Item(17, 2).Value = 0.168714240193367
$ws.Cells.Item(17, 3).Value = 0.22974132001399994
$ws.Cells.Item(17, 4).Value = 0.56414449214935303
$ws.Cells.Item(17, 5).Value = 0.53470790386199951
$ws.Cells.Item(18, 1).Value = 44317
$ws.Cells.Item(18, 2).Value = 0.020063307136297226
$ws.Cells.Item(18, 3).Value = -0.020555285736918449
$ws.Cells.Item(18, 4).Value = 0.54223060607910156
$ws.Cells.Item(18, 5).Value = 0.5078432559967041
$ws.Cells.Item(19, 1).Value = 44348
$ws.Cells.Item(19, 2).Value = 0.98783475160598755
$ws.Cells.Item(19, 3).Value = 0.96644324064254761
$ws.Cells.Item(19, 4).Value = 0.58736658096313477
$ws.Cells.Item(19, 5).Value = 0.55823147296905518
$ws.Cells.Item(20, 1).Value = 44378
$ws.Cells.Item(20, 2).Value = 0.86229372024536133
$ws.Cells.Item(20, 3).Value = 0.80186676979064941
$ws.Cells.Item(20, 4).Value = 0.58975428342819214
$ws.Cells.Item(20, 5).Value = 0.55724692344665527
$ws.Cells.Item(21, 1).Value = 44409
$ws.Cells.Item(21, 2).Value = 0.7834208607673645
$ws.Cells.Item(21, 3).Value = 0.70771795511245728
$ws.Cells.Item(21, 4).Value = 0.62779784202575684
$ws.Cells.Item(21, 5).Value = 0.59112966060638428
$ws.Cells.Item(22, 1).Value = 44440
$ws.Cells.Item(22, 2).Value = 0.56541162729263306
$ws.Cells.Item(22, 3).Value = 0.50165247917175293
$ws.Cells.Item(22, 4).Value = 0.67616474628448486
$ws.Cells.Item(22, 5).Value = 0.62921631336212158
$ws.Cells.Item(23, 1).Value = 44470
$ws.Cells.Item(23, 2).Value = 0.92950445413589478
$ws.Cells.Item(23, 3).Value = 0.90704536437988281
$ws.Cells.Item(23, 4).Value = 0.7726747989654541
$ws.Cells.Item(23, 5).Value = 0.7279888391494751
$ws.Cells.Item(24, 1).Value = 44501
$ws.Cells.Item(24, 2).Value = 0.75843614339828491
$ws.Cells.Item(24, 3).Value = 0.72591394186019897
$ws.Cells.Item(24, 4).Value = 0.72733175754547119
$ws.Cells.Item(24, 5).Value = 0.6796460747718811
$ws.Cells.Item(25, 1).Value = 44531
$ws.Cells.Item(25, 2).Value = 0.5745013952255249
$ws.Cells.Item(25, 3).Value = 0.50034099817276001
$ws.Cells.Item(25, 4).Value = 0.64191234111785889
$ws.Cells.Item(25, 5).Value = 0.58127343654632568
$ws.Cells.Item(26, 1).Value = 44562
$ws.Cells.Item(26, 2).Value = 0.60401636362075806
$ws.Cells.Item(26, 3).Value = 0.57252109050750732
$ws.Cells.Item(26, 4).Value = 0.54834198951721191
$ws.Cells.Item(26, 5).Value = 0.47506833076477051
$ws.Cells.Item(27, 1).Value = 44593
$ws.Cells.Item(27, 2).Value = 0.88865375518798828
$ws.Cells.Item(27, 3).Value = 0.86839777231216431
$ws.Cells.Item(27, 4).Value = 0.52398020029067993
$ws.Cells.Item(27, 5).Value = 0.44424682855606079
$ws.Cells.Item(28, 1).Value = 44621
$ws.Cells.Item(28, 2).Value = 0.57974767684936523
$ws.Cells.Item(28, 3).Value = 0.53135812282562256
$ws.Cells.Item(28, 4).Value = 0.44605356454849243
$ws.Cells.Item(28, 5).Value = 0.35953664779663086
$ws.Cells.Item(29, 1).Value = 44652
$ws.Cells.Item(29, 2).Value = 0.093518778681755066
$ws.Cells.Item(29, 3).Value = -0.083486981689929962
$ws.Cells.Item(29, 4).Value = 0.34521672129631042
$ws.Cells.Item(29, 5).Value = 0.24805796146392822
$ws.Cells.Item(30, 1).Value = 44682
$ws.Cells.Item(30, 2).Value = -0.058712445199489594
$ws.Cells.Item(30, 3).Value = -0.24812793731689453
$ws.Cells.Item(30, 4).Value = 0.27751493453979492
$ws.Cells.Item(30, 5).Value = 0.17702849209308624
$ws.Cells.Item(31, 1).Value = 44713
$ws.Cells.Item(31, 2).Value = 0.34615576267242432
$ws.Cells.Item(31, 3).Value = 0.22425921261310577
$ws.Cells.Item(31, 4).Value = 0.21286661922931671
$ws.Cells.Item(31, 5).Value = 0.11276738345623016
$ws.Cells.Item(32, 1).Value = 44743
$ws.Cells.Item(32, 2).Value = 0.22816456854343414
$ws.Cells.Item(32, 3).Value = 0.14465358853340149
$ws.Cells.Item(32, 4).Value = 0.11221088469028473
$ws.Cells.Item(32, 5).Value = 0.0037512965500354767
$ws.Cells.Item(33, 1).Value = 44774
$ws.Cells.Item(33, 2).Value = -0.14909525215625763
$ws.Cells.Item(33, 3).Value = -0.27739426493644714
$ws.Cells.Item(33, 4).Value = 0.029440870508551598
$ws.Cells.Item(33, 5).Value = -0.080971375107765198
$ws.Cells.Item(34, 1).Value = 44805
$ws.Cells.Item(34, 2).Value = -0.034814868122339249
$ws.Cells.Item(34, 3).Value = -0.13892413675785065
$ws.Cells.Item(34, 4).Value = 0.03338068351149559
$ws.Cells.Item(34, 5).Value = -0.071730323135852814
$ws.Cells.Item(35, 1).Value = 44835
$ws.Cells.Item(35, 2).Value = 0.022181600332260132
$ws.Cells.Item(35, 3).Value = -0.0058289356529712677
$ws.Cells.Item(35, 4).Value = 0.025342671200633049
$ws.Cells.Item(35, 5).Value = -0.083100423216819763
$ws.Cells.Item(36, 1).Value = 44866
$ws.Cells.Item(36, 2).Value = -0.017247855663299561
$ws.Cells.Item(36, 3).Value = -0.11274699866771698
$ws.Cells.Item(36, 4).Value = -0.023631807416677475
$ws.Cells.Item(36, 5).Value = -0.13116812705993652
$ws.Cells.Item(37, 1).Value = 44896
$ws.Cells.Item(37, 2).Value = -0.1651824563741684
$ws.Cells.Item(37, 3).Value = -0.2311459481716156
$ws.Cells.Item(37, 4).Value = -0.041741345077753067
$ws.Cells.Item(37, 5).Value = -0.16337338089942932
$ws.Cells.Item(38, 1).Value = 44927
$ws.Cells.Item(38, 2).Value = 0.12897710502147675
$ws.Cells.Item(38, 3).Value = -0.00031748833134770393
$ws.Cells.Item(38, 4).Value = -0.0066232206299901009
$ws.Cells.Item(38, 5).Value = -0.13906648755073547
$ws.Cells.Item(39, 1).Value = 44958
$ws.Cells.Item(39, 2).Value = -0.13105456531047821
$ws.Cells.Item(39, 3).Value = -0.35045883059501648
$ws.Cells.Item(39, 4).Value = 0.0067126024514436722
$ws.Cells.Item(39, 5).Value = -0.13146647810935974
$ws.Cells.Item(40, 1).Value = 44986
$ws.Cells.Item(40, 2).Value = -0.09461454302072525
$ws.Cells.Item(40, 3).Value = -0.20835006237030029
$ws.Cells.Item(40, 4).Value = 0.03950674831867218
$ws.Cells.Item(40, 5).Value = -0.11808700114488602
$ws.Cells.Item(41, 1).Value = 45017
$ws.Cells.Item(41, 2).Value = 0.065178714692592621
$ws.Cells.Item(41, 3).Value = -0.14519371092319489
$ws.Cells.Item(41, 4).Value = 0.076159395277500153
$ws.Cells.Item(41, 5).Value = -0.086544610559940338
$ws.Cells.Item(42, 1).Value = 45047
$ws.Cells.Item(42, 2).Value = 0.16696788370609283
$ws.Cells.Item(42, 3).Value = -0.058632321655750275
$ws.Cells.Item(42, 4).Value = 0.13985089957714081
$ws.Cells.Item(42, 5).Value = -0.027729466557502747
$ws.Cells.Item(43, 1).Value = 45078
$ws.Cells.Item(43, 2).Value = 0.085207536816596985
$ws.Cells.Item(43, 3).Value = -0.070524029433727264
$ws.Cells.Item(43, 4).Value = 0.17068886756896973
$ws.Cells.Item(43, 5).Value = 0.0063531943596899509
$ws.Cells.Item(44, 1).Value = 45108
$ws.Cells.Item(44, 2).Value = 0.31732892990112305
$ws.Cells.Item(44, 3).Value = 0.11458641290664673
$ws.Cells.Item(44, 4).Value = 0.23953601717948914
$ws.Cells.Item(44, 5).Value = 0.082219205796718597
$ws.Cells.Item(45, 1).Value = 45139
$ws.Cells.Item(45, 2).Value = 0.3126259446144104
$ws.Cells.Item(45, 3).Value = 0.17113450169563293
$ws.Cells.Item(45, 4).Value = 0.32293534278869629
$ws.Cells.Item(45, 5).Value = 0.16530118882656097
$ws.Cells.Item(46, 1).Value = 45170
$ws.Cells.Item(46, 2).Value = 0.40804114937782288
$ws.Cells.Item(46, 3).Value = 0.29819032549858093
$ws.Cells.Item(46, 4).Value = 0.37612941861152649
$ws.Cells.Item(46, 5).Value = 0.22521060705184937
$ws.Cells.Item(47, 1).Value = 45200
$ws.Cells.Item(47, 2).Value = 0.40651878714561462
$ws.Cells.Item(47, 3).Value = 0.30642646551132202
$ws.Cells.Item(47, 4).Value = 0.42497488856315613
$ws.Cells.Item(47, 5).Value = 0.2819570004940033
$ws.Cells.Item(48, 1).Value = 45231
$ws.Cells.Item(48, 2).Value = 0.48856979608535767
$ws.Cells.Item(48, 3).Value = 0.33233526349067688
$ws.Cells.Item(48, 4).Value = 0.47562271356582642
$ws.Cells.Item(48, 5).Value = 0.32222369313240051
$ws.Cells.Item(49, 1).Value = 45261
$ws.Cells.Item(49, 2).Value = 0.65597927570343018
$ws.Cells.Item(49, 3).Value = 0.53938776254653931
$ws.Cells.Item(49, 4).Value = 0.50628876686096191
$ws.Cells.Item(49, 5).Value = 0.35612225532531738
$ws.Cells.Item(50, 1).Value = 45292
$ws.Cells.Item(50, 2).Value = 0.54392540454864502
$ws.Cells.Item(50, 3).Value = 0.39399105310440063
$ws.Cells.Item(50, 4).Value = 0.53772157430648804
$ws.Cells.Item(50, 5).Value = 0.37842926383018494
$ws.Cells.Item(51, 1).Value = 45323
$ws.Cells.Item(51, 2).Value = 0.60657721757888794
$ws.Cells.Item(51, 3).Value = 0.45208531618118286
$ws.Cells.Item(51, 4).Value = 0.54334825277328491
$ws.Cells.Item(51, 5).Value = 0.37346422672271729
$ws.Cells.Item(52, 1).Value = 45352
$ws.Cells.Item(52, 2).Value = 0.54103797674179077
$ws.Cells.Item(52, 3).Value = 0.2918761670589447
$ws.Cells.Item(52, 4).Value = 0.58508741855621338
$ws.Cells.Item(52, 5).Value = 0.40517356991767883
$ws.Cells.Item(53, 1).Value = 45383
$ws.Cells.Item(53, 2).Value = 0.5933234691619873
$ws.Cells.Item(53, 3).Value = 0.41967332363128662
$ws.Cells.Item(53, 4).Value = 0.64761435985565186
$ws.Cells.Item(53, 5).Value = 0.46538490056991577
$ws.Cells.Item(54, 1).Value = 45413
$ws.Cells.Item(54, 2).Value = 0.59552127122879028
$ws.Cells.Item(54, 3).Value = 0.37189778685569763
$ws.Cells.Item(54, 4).Value = 0.68400126695632935
$ws.Cells.Item(54, 5).Value = 0.50265538692474365
$ws.Cells.Item(55, 1).Value = 45444
$ws.Cells.Item(55, 2).Value = 0.45868101716041565
$ws.Cells.Item(55, 3).Value = 0.25350496172904968
$ws.Cells.Item(55, 4).Value = 0.70473229885101318
$ws.Cells.Item(55, 5).Value = 0.52199298143386841
$ws.Cells.Item(56, 1).Value = 45474
$ws.Cells.Item(56, 2).Value = 0.78217107057571411
$ws.Cells.Item(56, 3).Value = 0.59181046485900879
$ws.Cells.Item(56, 4).Value = 0.73539263010025024
$ws.Cells.Item(56, 5).Value = 0.5561710000038147
$ws.Cells.Item(57, 1).Value = 45505
$ws.Cells.Item(57, 2).Value = 1.0513128042221069
$ws.Cells.Item(57, 3).Value = 0.87423735857009888
$ws.Cells.Item(57, 4).Value = 0.75555700063705444
$ws.Cells.Item(57, 5).Value = 0.59127509593963623
$ws.Cells.Item(58, 1).Value = 45536
$ws.Cells.Item(58, 2).Value = 0.98346120119094849
$ws.Cells.Item(58, 3).Value = 0.87482202053070068
$ws.Cells.Item(58, 4).Value = 0.79562336206436157
$ws.Cells.Item(58, 5).Value = 0.62913239002227783
$ws.Cells.Item(59, 1).Value = 45566
$ws.Cells.Item(59, 2).Value = 0.73050457239151001
$ws.Cells.Item(59, 3).Value = 0.56802940368652344
$ws.Cells.Item(59, 4).Value = 0.80862295627593994
$ws.Cells.Item(59, 5).Value = 0.63972145318984985
$ws.Cells.Item(60, 1).Value = 45597
$ws.Cells.Item(60, 2).Value = 0.88252043724060059
$ws.Cells.Item(60, 3).Value = 0.75968754291534424
$ws.Cells.Item(60, 4).Value = 0.83484983444213867
$ws.Cells.Item(60, 5).Value = 0.67448800802230835
$ws.Cells.Item(61, 1).Value = 45627
$ws.Cells.Item(61, 2).Value = 0.72251707315444946
$ws.Cells.Item(61, 3).Value = 0.60781276226043701
$ws.Cells.Item(61, 4).Value = 0.84143465757369995
$ws.Cells.Item(61, 5).Value = 0.684822678565979
$ws.Cells.Item(62, 1).Value = 45658
$ws.Cells.Item(62, 2).Value = 0.95392090082168579
$ws.Cells.Item(62, 3).Value = 0.76038902997970581
$ws.Cells.Item(62, 4).Value = 0.81145209074020386
$ws.Cells.Item(62, 5).Value = 0.65776348114013672
$ws.Cells.Item(63, 1).Value = 45689
$ws.Cells.Item(63, 2).Value = 0.71251761913299561
$ws.Cells.Item(63, 3).Value = 0.4671996533870697
$ws.Cells.Item(63, 4).Value = 0.7827838659286499
$ws.Cells.Item(63, 5).Value = 0.62158703804016113
$ws.Cells.Item(64, 1).Value = 45717
$ws.Cells.Item(64, 2).Value = 0.69472277164459229
$ws.Cells.Item(64, 3).Value = 0.56640380620956421
$ws.Cells.Item(64, 4).Value = 0.7932397723197937
$ws.Cells.Item(64, 5).Value = 0.63229858875274658
